# -----------------------------------------------------------------------
# Weekly refresh for "Vega Modelo de Temuco - Coco": a new week of price
# observations is inserted at the front of the tracked date window, which
# shifts every later rows Fecha/Volumen/Precio one slot down the series.
# The two oldest rows that fall out of the window are appended at the end
# (rows 58-59), and the sheets used range grows from T57 to T59.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = 44664 ; $ws.Range("N28").Value = 28000 ; $ws.Range("O28").Value = 28000 ; $ws.Range("P28").Value = 28000 ; $ws.Range("S28").Value = 1400
$ws.Range("D29").Value = 44175 ; $ws.Range("M29").Value = 25 ; $ws.Range("N29").Value = 23000 ; $ws.Range("O29").Value = 23000 ; $ws.Range("P29").Value = 23000 ; $ws.Range("S29").Value = 1150
$ws.Range("D30").Value = 44488 ; $ws.Range("M30").Value = 40 ; $ws.Range("N30").Value = 20000 ; $ws.Range("O30").Value = 20000 ; $ws.Range("P30").Value = 20000 ; $ws.Range("S30").Value = 1000
$ws.Range("D31").Value = 44222 ; $ws.Range("M31").Value = 15 ; $ws.Range("N31").Value = 25000 ; $ws.Range("O31").Value = 25000 ; $ws.Range("P31").Value = 25000 ; $ws.Range("S31").Value = 1250
$ws.Range("D32").Value = 44421 ; $ws.Range("M32").Value = 20
$ws.Range("D33").Value = 44469 ; $ws.Range("M33").Value = 40 ; $ws.Range("N33").Value = 24000 ; $ws.Range("O33").Value = 24000 ; $ws.Range("P33").Value = 24000 ; $ws.Range("S33").Value = 1200
$ws.Range("D34").Value = 44334 ; $ws.Range("N34").Value = 25000 ; $ws.Range("O34").Value = 25000 ; $ws.Range("P34").Value = 25000 ; $ws.Range("S34").Value = 1250
$ws.Range("D35").Value = 44475 ; $ws.Range("M35").Value = 20
$ws.Range("D36").Value = 44356 ; $ws.Range("M36").Value = 15
$ws.Range("D37").Value = 44363 ; $ws.Range("M37").Value = 30
$ws.Range("D38").Value = 44424 ; $ws.Range("N38").Value = 24000 ; $ws.Range("O38").Value = 24000 ; $ws.Range("P38").Value = 24000 ; $ws.Range("S38").Value = 1200
$ws.Range("D39").Value = 44442 ; $ws.Range("M39").Value = 25 ; $ws.Range("N39").Value = 23000 ; $ws.Range("O39").Value = 23000 ; $ws.Range("P39").Value = 23000 ; $ws.Range("S39").Value = 1150
$ws.Range("D40").Value = 44462 ; $ws.Range("M40").Value = 10 ; $ws.Range("N40").Value = 24000 ; $ws.Range("O40").Value = 24000 ; $ws.Range("P40").Value = 24000 ; $ws.Range("S40").Value = 1200
$ws.Range("D41").Value = 44214 ; $ws.Range("M41").Value = 15 ; $ws.Range("N41").Value = 25000 ; $ws.Range("O41").Value = 25000 ; $ws.Range("P41").Value = 25000 ; $ws.Range("S41").Value = 1250
$ws.Range("D42").Value = 44645 ; $ws.Range("M42").Value = 10 ; $ws.Range("N42").Value = 28000 ; $ws.Range("O42").Value = 28000 ; $ws.Range("P42").Value = 28000 ; $ws.Range("S42").Value = 1400
$ws.Range("D43").Value = 44235 ; $ws.Range("M43").Value = 15
$ws.Range("D44").Value = 44412 ; $ws.Range("N44").Value = 25000 ; $ws.Range("O44").Value = 25000 ; $ws.Range("P44").Value = 25000 ; $ws.Range("S44").Value = 1250
$ws.Range("D45").Value = 44468 ; $ws.Range("M45").Value = 20 ; $ws.Range("N45").Value = 24000 ; $ws.Range("O45").Value = 24000 ; $ws.Range("P45").Value = 24000 ; $ws.Range("S45").Value = 1200
$ws.Range("D46").Value = 44663 ; $ws.Range("M46").Value = 20 ; $ws.Range("N46").Value = 28000 ; $ws.Range("O46").Value = 28000 ; $ws.Range("P46").Value = 28000 ; $ws.Range("S46").Value = 1400
$ws.Range("D47").Value = 44251 ; $ws.Range("M47").Value = 15 ; $ws.Range("N47").Value = 25000 ; $ws.Range("O47").Value = 25000 ; $ws.Range("P47").Value = 25000 ; $ws.Range("S47").Value = 1250
$ws.Range("D48").Value = 44452 ; $ws.Range("M48").Value = 25 ; $ws.Range("N48").Value = 25000 ; $ws.Range("O48").Value = 25000 ; $ws.Range("P48").Value = 25000 ; $ws.Range("S48").Value = 1250
$ws.Range("D49").Value = 44434 ; $ws.Range("M49").Value = 20 ; $ws.Range("N49").Value = 24000 ; $ws.Range("O49").Value = 24000 ; $ws.Range("P49").Value = 24000 ; $ws.Range("S49").Value = 1200
$ws.Range("D50").Value = 44432 ; $ws.Range("M50").Value = 30
$ws.Range("D51").Value = 44231 ; $ws.Range("M51").Value = 15 ; $ws.Range("N51").Value = 25000 ; $ws.Range("O51").Value = 25000 ; $ws.Range("P51").Value = 25000 ; $ws.Range("S51").Value = 1250
$ws.Range("D52").Value = 44428
$ws.Range("D53").Value = 44389 ; $ws.Range("M53").Value = 20
$ws.Range("D54").Value = 44426 ; $ws.Range("M54").Value = 15
$ws.Range("D55").Value = 44396 ; $ws.Range("M55").Value = 12
$ws.Range("D56").Value = 44435 ; $ws.Range("M56").Value = 100
$ws.Range("D57").Value = 44467 ; $ws.Range("M57").Value = 20
# Two brand-new rows appended at the bottom (rows 58 and 59) carrying the
# oldest two observations that rolled out of the front of the window;
# every non-numeric column below mirrors the constant market/product
# attributes already used throughout this sheet (rows 2-57).

# Row 58
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44349
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100108
$ws.Range("H58").Value = "Tropicales y subtropicales"
$ws.Range("I58").Value = 100108007
$ws.Range("J58").Value = "Coco"
$ws.Range("K58").Value = "Sin especificar"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 30
$ws.Range("N58").Value = 24000
$ws.Range("O58").Value = 24000
$ws.Range("P58").Value = 24000
$ws.Range("Q58").Value = "$/malla 20 unidades"
$ws.Range("R58").Value = "Perú"
$ws.Range("S58").Value = 1200
$ws.Range("T58").Value = 20
$ws.Range("D58").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 59
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 44466
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100108
$ws.Range("H59").Value = "Tropicales y subtropicales"
$ws.Range("I59").Value = 100108007
$ws.Range("J59").Value = "Coco"
$ws.Range("K59").Value = "Sin especificar"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 70
$ws.Range("N59").Value = 24000
$ws.Range("O59").Value = 24000
$ws.Range("P59").Value = 24000
$ws.Range("Q59").Value = "$/malla 20 unidades"
$ws.Range("R59").Value = "Perú"
$ws.Range("S59").Value = 1200
$ws.Range("T59").Value = 20
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"